$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update total_customers (C), returning_customers (D), new_customers (E), recurrence_rate (F)
# for rows 2-28 as per updated BIBI data

$ws.Cells.Item(2, 3).Value = 513
$ws.Cells.Item(2, 4).Value = 35
$ws.Cells.Item(2, 5).Value = 478
$ws.Cells.Item(2, 6).Value = 6.849315068493151

$ws.Cells.Item(3, 3).Value = 398
$ws.Cells.Item(3, 4).Value = 28
$ws.Cells.Item(3, 5).Value = 370
$ws.Cells.Item(3, 6).Value = 5.458089668615984

$ws.Cells.Item(4, 3).Value = 316
$ws.Cells.Item(4, 4).Value = 16
$ws.Cells.Item(4, 5).Value = 300
$ws.Cells.Item(4, 6).Value = 4.020100502512562

$ws.Cells.Item(5, 3).Value = 563
$ws.Cells.Item(5, 4).Value = 20
$ws.Cells.Item(5, 5).Value = 543
$ws.Cells.Item(5, 6).Value = 6.329113924050633

$ws.Cells.Item(6, 3).Value = 1025
$ws.Cells.Item(6, 4).Value = 22
$ws.Cells.Item(6, 5).Value = 1003
$ws.Cells.Item(6, 6).Value = 3.907637655417407

$ws.Cells.Item(7, 3).Value = 723
$ws.Cells.Item(7, 4).Value = 50
$ws.Cells.Item(7, 5).Value = 673
$ws.Cells.Item(7, 6).Value = 4.878048780487805

$ws.Cells.Item(8, 3).Value = 656
$ws.Cells.Item(8, 4).Value = 20
$ws.Cells.Item(8, 5).Value = 636
$ws.Cells.Item(8, 6).Value = 2.766251728907331

$ws.Cells.Item(9, 3).Value = 897
$ws.Cells.Item(9, 4).Value = 15
$ws.Cells.Item(9, 5).Value = 882
$ws.Cells.Item(9, 6).Value = 2.286585365853659

$ws.Cells.Item(10, 3).Value = 666
$ws.Cells.Item(10, 4).Value = 48
$ws.Cells.Item(10, 5).Value = 618
$ws.Cells.Item(10, 6).Value = 5.351170568561873

$ws.Cells.Item(11, 3).Value = 515
$ws.Cells.Item(11, 4).Value = 17
$ws.Cells.Item(11, 5).Value = 498
$ws.Cells.Item(11, 6).Value = 2.552552552552553

$ws.Cells.Item(12, 3).Value = 696
$ws.Cells.Item(12, 4).Value = 35
$ws.Cells.Item(12, 5).Value = 661
$ws.Cells.Item(12, 6).Value = 6.796116504854369

$ws.Cells.Item(13, 3).Value = 689
$ws.Cells.Item(13, 4).Value = 38
$ws.Cells.Item(13, 5).Value = 651
$ws.Cells.Item(13, 6).Value = 5.459770114942529

$ws.Cells.Item(14, 3).Value = 744
$ws.Cells.Item(14, 4).Value = 45
$ws.Cells.Item(14, 5).Value = 699
$ws.Cells.Item(14, 6).Value = 6.531204644412192

$ws.Cells.Item(15, 3).Value = 730
$ws.Cells.Item(15, 4).Value = 49
$ws.Cells.Item(15, 5).Value = 681
$ws.Cells.Item(15, 6).Value = 6.586021505376344

$ws.Cells.Item(16, 3).Value = 748
$ws.Cells.Item(16, 4).Value = 66
$ws.Cells.Item(16, 5).Value = 682
$ws.Cells.Item(16, 6).Value = 9.04109589041096

$ws.Cells.Item(17, 3).Value = 820
$ws.Cells.Item(17, 4).Value = 77
$ws.Cells.Item(17, 5).Value = 743
$ws.Cells.Item(17, 6).Value = 10.29411764705882

$ws.Cells.Item(18, 3).Value = 723
$ws.Cells.Item(18, 4).Value = 58
$ws.Cells.Item(18, 5).Value = 665
$ws.Cells.Item(18, 6).Value = 7.073170731707316

$ws.Cells.Item(19, 3).Value = 795
$ws.Cells.Item(19, 4).Value = 48
$ws.Cells.Item(19, 5).Value = 747
$ws.Cells.Item(19, 6).Value = 6.639004149377594

$ws.Cells.Item(20, 3).Value = 762
$ws.Cells.Item(20, 4).Value = 58
$ws.Cells.Item(20, 5).Value = 704
$ws.Cells.Item(20, 6).Value = 7.29559748427673

$ws.Cells.Item(21, 3).Value = 743
$ws.Cells.Item(21, 4).Value = 50
$ws.Cells.Item(21, 5).Value = 693
$ws.Cells.Item(21, 6).Value = 6.561679790026247

$ws.Cells.Item(22, 3).Value = 669
$ws.Cells.Item(22, 4).Value = 38
$ws.Cells.Item(22, 5).Value = 631
$ws.Cells.Item(22, 6).Value = 5.114401076716016

$ws.Cells.Item(23, 3).Value = 651
$ws.Cells.Item(23, 4).Value = 49
$ws.Cells.Item(23, 5).Value = 602
$ws.Cells.Item(23, 6).Value = 7.324364723467862

$ws.Cells.Item(24, 3).Value = 594
$ws.Cells.Item(24, 4).Value = 62
$ws.Cells.Item(24, 5).Value = 532
$ws.Cells.Item(24, 6).Value = 9.523809523809524

$ws.Cells.Item(25, 3).Value = 711
$ws.Cells.Item(25, 4).Value = 63
$ws.Cells.Item(25, 5).Value = 648
$ws.Cells.Item(25, 6).Value = 10.60606060606061

$ws.Cells.Item(26, 3).Value = 737
$ws.Cells.Item(26, 4).Value = 48
$ws.Cells.Item(26, 5).Value = 689
$ws.Cells.Item(26, 6).Value = 6.751054852320674

$ws.Cells.Item(27, 3).Value = 642
$ws.Cells.Item(27, 4).Value = 63
$ws.Cells.Item(27, 5).Value = 579
$ws.Cells.Item(27, 6).Value = 8.548168249660787

$ws.Cells.Item(28, 3).Value = 142
$ws.Cells.Item(28, 4).Value = 23
$ws.Cells.Item(28, 5).Value = 119
$ws.Cells.Item(28, 6).Value = 3.582554517133956
